# Apply the "Big changes" commit: append a new trade row (row 8) to the
# HZNP trade log sheet and widen column C to fit the new, wider value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 7) down into
# the new row 8, so the date/boolean cells keep the same number formats
# and styles (s="1") as the rest of the table.
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new trade's data.
$ws.Cells.Item(8, 1).Value = 42649.644837962966   # Date
$ws.Cells.Item(8, 2).Value = $true                # Profitable
$ws.Cells.Item(8, 3).Value = 10068.83             # Principle
$ws.Cells.Item(8, 4).Value = 9957.7999999999993   # Start Principle
$ws.Cells.Item(8, 5).Value = 18.829999999999998   # BuyPrice
$ws.Cells.Item(8, 6).Value = 19.25                # SellPrice
$ws.Cells.Item(8, 7).Value = $false               # IsShortSell
$ws.Cells.Item(8, 8).Value = 2.23                 # Price Change %
$ws.Cells.Item(8, 9).Value = $false               # Strong trade

# The new Principle value (10068.83) is wider than the previous longest
# entry in column C, so the sheet's best-fit column width grows.
$ws.Columns.Item(3).ColumnWidth = 8
